$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 462 (shifts the existing rows 462..479 down to 463..480)
$ws.Rows.Item(462).Insert()

# Populate the newly inserted row 462 with the new weekly record
$ws.Range("A462").Value = 7
$ws.Range("B462").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C462").Value = "Ñuble"
$ws.Range("D462").Value = 44939
$ws.Range("E462").Value = 16
$ws.Range("F462").Value = 100114001
$ws.Range("G462").Value = "Papa"
$ws.Range("H462").Value = "Asterix"
$ws.Range("I462").Value = "1a (cosecha)"
$ws.Range("J462").Value = 200
$ws.Range("K462").Value = 11000
$ws.Range("L462").Value = 12000
$ws.Range("M462").Value = 11500
$ws.Range("N462").Value = "`$/saco 25 kilos"
$ws.Range("O462").Value = "Región del Maule"
$ws.Range("P462").Value = 460
$ws.Range("Q462").Value = 25
$ws.Range("R462").Value = "Hortaliza"
